{"js": "// TECHNICAL SKILL bullet list updates (resume edited to match IELTS-general refresh).\nconst body = context.document.body;\n\n// 1) \"Used efficient programming languages as: Java, C#, python\"\n//    -> \"Used efficient programming languages as: Java, C#, Python, ASP.NET.\"\nlet results = body.search(\"Used efficient programming languages as: Java, C#, python\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Used efficient programming languages as: Java, C#, Python, ASP.NET.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) \"Framework for developing mobile applications as React Native CLI or Expo\"\n//    -> same text with a trailing period added.\nresults = body.search(\"Framework for developing mobile applications as React Native CLI or Expo\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Framework for developing mobile applications as React Native CLI or Expo.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 3) \"...styles using TailwindCSS\" -> add trailing period.\nresults = body.search(\"TailwindCSS\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\".\", \"After\");\n  await context.sync();\n}\n\n// 4) \"Backend using Express with Sqlite database, MongoDB.\"\n//    -> \"Backend using Express with SQLlite database, MongoDB, Oracle, SQLServer.\"\nresults = body.search(\"Backend using Express with Sqlite database, MongoDB.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Backend using Express with SQLlite database, MongoDB, Oracle, SQLServer.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# TECHNICAL SKILL bullet list updates (resume edited to match IELTS-general refresh).\n$d = $word.ActiveDocument\n\n# 1) \"Used efficient programming languages as: Java, C#, python\"\n#    -> \"Used efficient programming languages as: Java, C#, Python, ASP.NET.\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"Used efficient programming languages as: Java, C#, python\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Used efficient programming languages as: Java, C#, Python, ASP.NET.\",\n    2\n)\n\n# 2) \"Framework for developing mobile applications as React Native CLI or Expo\"\n#    -> same text with a trailing period added.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"Framework for developing mobile applications as React Native CLI or Expo\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Framework for developing mobile applications as React Native CLI or Expo.\",\n    2\n)\n\n# 3) \"...styles using TailwindCSS\" -> add trailing period.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"TailwindCSS\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"TailwindCSS.\",\n    2\n)\n\n# 4) \"Backend using Express with Sqlite database, MongoDB.\"\n#    -> \"Backend using Express with SQLlite database, MongoDB, Oracle, SQLServer.\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"Backend using Express with Sqlite database, MongoDB.\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Backend using Express with SQLlite database, MongoDB, Oracle, SQLServer.\",\n    2\n)\n"}
